# Apply WRI input-data update to the GDPbES (Guaranteed Dispatch Percentage by
# Electricity Source) sheet:
#   - Add a column header / units note in A1: "Guaranteed Dispatch Fraction
#     (dimensionless)", bold + wrap-text, and grow row 1's height to fit it.
#   - Add three new guaranteed-dispatch-source rows (crude oil, heavy or
#     residual fuel oil, municipal solid waste) that mirror existing rows
#     (petroleum for the two oil rows, biomass for municipal solid waste).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GDPbES")

# --- Row 15: crude oil (mirrors row 11, petroleum) -------------------------
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15").Formula = "=B11"
$ws.Range("C15:AK15").Formula = "=C11"

# --- Row 16: heavy or residual fuel oil (mirrors row 11, petroleum) -------
$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16").Formula = "=B11"
$ws.Range("C16:AK16").Formula = "=C11"

# --- Row 17: municipal solid waste (mirrors row 9, biomass) ---------------
$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17").Formula = "=B9"
$ws.Range("C17:AK17").Formula = "=C9"

# --- Header cell (A1): units label, bold + wrap, taller row ---------------
$ws.Range("A1").Value = "Guaranteed Dispatch Fraction (dimensionless)"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 45
